$wb = $excel.ActiveWorkbook

# --- Worksheet references ---
$wsLivestock = $wb.Worksheets.Item("livestock")
$wsExcretion = $wb.Worksheets.Item("excretion")
$wsProd      = $wb.Worksheets.Item("prod")
$wsGlobal    = $wb.Worksheets.Item("global")

# --- Rename the "Methanization power" header from tMB to tFW ---
# (affects the excretion sheet G1 header and the prod sheet H1 header,
#  both of which shared the same text before the edit)
$wsExcretion.Range("G1").Value = "Methanization power (MWh/tFW)"
$wsProd.Range("H1").Value = "Methanization power (MWh/tFW)"

# --- Update methanizer figures on the "global" sheet ---
$wsGlobal.Range("B6").Value = 1000
$wsGlobal.Range("B9").Value = 50000

# --- Add a new "Weight import" parameter row at the bottom of "global" ---
$wsGlobal.Range("A10").Value = "Weight import"
$wsGlobal.Range("B10").Value = 0

# --- Update cell selections left behind by the editing session ---
$wsLivestock.Activate()
$wsLivestock.Range("B7").Select()

$wsExcretion.Activate()
$wsExcretion.Range("G2").Select()

$wsProd.Activate()
$wsProd.Range("H1").Select()

$wsGlobal.Activate()
$wsGlobal.Range("B10").Select()
